$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gStation")

# winch.dwinch_dt (A1/B1): 100 -> 90 ("USe D by d of 90")
$ws.Range("B1").Value = 90

# winch.material (A8/B8): 1 -> 2 ("Use steel for winch")
$ws.Range("B8").Value = 2

# Update the active selection to B1 on this sheet
$ws.Activate()
$ws.Range("B1").Select()
